$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2107438016528926
$ws.Range("C2").Value = 0.5289256198347108
$ws.Range("J2").Value = 0.004132231404958678
$ws.Range("P2").Value = 0.1570247933884298
$ws.Range("S2").Value = 0.09917355371900827
$ws.Range("C3").Value = 0.01538461538461539
$ws.Range("J3").Value = 0.03846153846153846
$ws.Range("P3").Value = 0.7461538461538462
$ws.Range("S3").Value = 0.2
$ws.Range("O4").Value = 0.02857142857142857
$ws.Range("P4").Value = 0.8285714285714286
$ws.Range("S4").Value = 0.1428571428571428
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.5
$ws.Range("B6").Value = 0.08152173913043478
$ws.Range("D6").Value = 0.005434782608695652
$ws.Range("E6").Value = 0.005434782608695652
$ws.Range("F6").Value = 0.07608695652173914
$ws.Range("J6").Value = 0.2717391304347826
$ws.Range("O6").Value = 0.02717391304347826
$ws.Range("Q6").Value = 0.108695652173913
$ws.Range("R6").Value = 0.03260869565217391
$ws.Range("S6").Value = 0.391304347826087
$ws.Range("B7").Value = 0.09782608695652174
$ws.Range("D7").Value = 0.0108695652173913
$ws.Range("F7").Value = 0.07065217391304347
$ws.Range("J7").Value = 0.09239130434782608
$ws.Range("O7").Value = 0.0108695652173913
$ws.Range("Q7").Value = 0.2010869565217391
$ws.Range("R7").Value = 0.07065217391304347
$ws.Range("S7").Value = 0.4456521739130435
$ws.Range("B8").Value = 0.1004784688995215
$ws.Range("D8").Value = 0.01913875598086124
$ws.Range("E8").Value = 0.004784688995215311
$ws.Range("F8").Value = 0.05980861244019139
$ws.Range("J8").Value = 0.09090909090909091
$ws.Range("O8").Value = 0.01913875598086124
$ws.Range("Q8").Value = 0.1722488038277512
$ws.Range("R8").Value = 0.09330143540669857
$ws.Range("S8").Value = 0.4401913875598086
$ws.Range("B9").Value = 0.1
$ws.Range("D9").Value = 0.005
$ws.Range("F9").Value = 0.05
$ws.Range("J9").Value = 0.11
$ws.Range("O9").Value = 0.025
$ws.Range("Q9").Value = 0.225
$ws.Range("R9").Value = 0.065
$ws.Range("S9").Value = 0.42
$ws.Range("B10").Value = 0.09313725490196079
$ws.Range("D10").Value = 0.02352941176470588
$ws.Range("E10").Value = 0.00196078431372549
$ws.Range("F10").Value = 0.0607843137254902
$ws.Range("J10").Value = 0.1049019607843137
$ws.Range("O10").Value = 0.0196078431372549
$ws.Range("Q10").Value = 0.2333333333333333
$ws.Range("R10").Value = 0.07058823529411765
$ws.Range("S10").Value = 0.392156862745098
$ws.Range("G11").Value = 0.1043478260869565
$ws.Range("J11").Value = 0.05217391304347826
$ws.Range("K11").Value = 0.1434782608695652
$ws.Range("L11").Value = 0.6869565217391305
$ws.Range("S11").Value = 0.01304347826086956
$ws.Range("G12").Value = 0.8395061728395061
$ws.Range("J12").Value = 0.1296296296296296
$ws.Range("L12").Value = 0.02469135802469136
$ws.Range("S12").Value = 0.006172839506172839
$ws.Range("G13").Value = 0.7906976744186046
$ws.Range("J13").Value = 0.186046511627907
$ws.Range("S13").Value = 0.02325581395348837
$ws.Range("F15").Value = 0.0179372197309417
$ws.Range("H15").Value = 0.1838565022421525
$ws.Range("I15").Value = 0.07174887892376682
$ws.Range("J15").Value = 0.304932735426009
$ws.Range("K15").Value = 0.07623318385650224
$ws.Range("M15").Value = 0.0179372197309417
$ws.Range("O15").Value = 0.08520179372197309
$ws.Range("S15").Value = 0.242152466367713
$ws.Range("F16").Value = 0.01219512195121951
$ws.Range("H16").Value = 0.2621951219512195
$ws.Range("I16").Value = 0.09146341463414634
$ws.Range("J16").Value = 0.3292682926829268
$ws.Range("K16").Value = 0.1097560975609756
$ws.Range("M16").Value = 0.03658536585365853
$ws.Range("O16").Value = 0.06097560975609756
$ws.Range("S16").Value = 0.0975609756097561
$ws.Range("F17").Value = 0.02941176470588235
$ws.Range("H17").Value = 0.1862745098039216
$ws.Range("I17").Value = 0.1225490196078431
$ws.Range("J17").Value = 0.3676470588235294
$ws.Range("K17").Value = 0.08823529411764706
$ws.Range("M17").Value = 0.01470588235294118
$ws.Range("O17").Value = 0.07352941176470588
$ws.Range("S17").Value = 0.1176470588235294
$ws.Range("F18").Value = 0.01418439716312057
$ws.Range("H18").Value = 0.148936170212766
$ws.Range("I18").Value = 0.0851063829787234
$ws.Range("J18").Value = 0.4609929078014184
$ws.Range("K18").Value = 0.07092198581560284
$ws.Range("M18").Value = 0.01418439716312057
$ws.Range("O18").Value = 0.1063829787234043
$ws.Range("S18").Value = 0.09929078014184398
$ws.Range("F19").Value = 0.01358695652173913
$ws.Range("H19").Value = 0.2210144927536232
$ws.Range("I19").Value = 0.09782608695652174
$ws.Range("J19").Value = 0.3704710144927536
$ws.Range("K19").Value = 0.1041666666666667
$ws.Range("M19").Value = 0.02445652173913044
$ws.Range("O19").Value = 0.07789855072463768
